# Insert two new data rows (old row 138/139 area) for Zapallo italiano,
# shifting the existing rows 138..193 down to 140..195, then populate the
# two newly-inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 138 - this shifts rows
# 138:193 down to 140:195 and extends the used range / dimension
# automatically.
$ws.Rows("138:139").Insert()

# New row 138
$ws.Cells.Item(138, 1).Value = 9
$ws.Cells.Item(138, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(138, 3).Value = "Metropolitana"
$ws.Cells.Item(138, 4).Value = 44455
$ws.Cells.Item(138, 5).Value = 13
$ws.Cells.Item(138, 6).Value = 100112032
$ws.Cells.Item(138, 7).Value = "Zapallo italiano"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 88
$ws.Cells.Item(138, 11).Value = 16000
$ws.Cells.Item(138, 12).Value = 17000
$ws.Cells.Item(138, 13).Value = 16500
$ws.Cells.Item(138, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(138, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(138, 16).Value = 275
$ws.Cells.Item(138, 17).Value = 60
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# New row 139
$ws.Cells.Item(139, 1).Value = 9
$ws.Cells.Item(139, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(139, 3).Value = "Metropolitana"
$ws.Cells.Item(139, 4).Value = 44455
$ws.Cells.Item(139, 5).Value = 13
$ws.Cells.Item(139, 6).Value = 100112032
$ws.Cells.Item(139, 7).Value = "Zapallo italiano"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Segunda"
$ws.Cells.Item(139, 10).Value = 61
$ws.Cells.Item(139, 11).Value = 14000
$ws.Cells.Item(139, 12).Value = 15000
$ws.Cells.Item(139, 13).Value = 14492
$ws.Cells.Item(139, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(139, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value = 145
$ws.Cells.Item(139, 17).Value = 100
$ws.Cells.Item(139, 18).Value = "Hortaliza"
